# SSU-Pravljenje Spila.docx edit script
# Restrikovano pravljenje spila na registrovane korisnike

$d = $word.ActiveDocument

function Split-RunsInRange($d, $absStart, $offsets) {
    # $offsets is a list of (start,end) pairs (relative to $absStart) that
    # should each become their own run. Applying a harmless formatting
    # toggle (Bold on/off) over exactly that span forces Word to give it
    # its own run without altering the visible formatting (since the net
    # effect of On-then-Off is a no-op).
    foreach ($pair in $offsets) {
        $s = $absStart + $pair[0]
        $e = $absStart + $pair[1]
        $r = $d.Range($s, $e)
        $r.Bold = 1
        $r.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 1) "Svako moze da napravi ..." paragraph -> restricted-to-logged-in text
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(111)
$rng1 = $p1.Range
$p1Start = $rng1.Start
$p1End = $rng1.End - 1
$target1 = $d.Range($p1Start, $p1End)
$newText1 = "Ulogovani korisnici mogu da naprave sopstveni špil za igru. Formiranje špila uključuje određivanje broja pojedinačnih karata u špilu kao i posebna pravila za specifične karte u špilu. Kad naprave špil mogu da započnu igru sa njim ili da ga sačuvaju na nalog."
$target1.Text = $newText1

$offsets1 = @(
    @(24,34),
    @(34,35),
    @(35,183),
    @(183,184),
    @(184,200),
    @(200,225),
    @(225,229),
    @(229,230),
    @(230,233)
)
Split-RunsInRange $d $p1Start $offsets1

# ---------------------------------------------------------------------
# 2) "biranjem istih iy liste" typo fix -> "iz"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(118)
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("biranjem istih iy liste", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2Start = $find2.Parent.Start
$p2End = $find2.Parent.End
$target2 = $d.Range($p2Start, $p2End)
$target2.Text = " biranjem istih iz liste"

$offsets2 = @(
    @(17,18)
)
Split-RunsInRange $d $p2Start $offsets2

# ---------------------------------------------------------------------
# 3) Remove the stray _GoBack bookmark left over in "1a Korisnik izadje..."
#    paragraph: an edit that spans across the bookmark's position causes
#    it to be dropped, matching the target (no text change there).
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(126)
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute("a Korisnik izađe", $true, $false, $false, $false, $false, $true, 1, $false, "a Korisnik izađe", 2)

# ---------------------------------------------------------------------
# 4) ".1 Ukoliko je korisnik ulogovan moze..." -> restricted rationale
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(132)
$find4 = $d.Content.Find
$find4.ClearFormatting()
$find4.Execute("U", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# Find the specific "U" that starts this sentence (after ".1 ") by scanning
# within paragraph 132's range instead of a global find.
$p4Start = $p4.Range.Start
$p4End = $p4.Range.End - 1
$p4Text = $p4.Range.Text
$uRelOffset = $p4Text.IndexOf("U")
$uAbs = $p4Start + $uRelOffset
$target4 = $d.Range($uAbs, $p4End)
$newText4 = "Umesto pokretanja igre, korisnik može i da sačuva špil (dugme " + [char]8217 + "Save" + [char]8217 + ")."
$target4.Text = $newText4

$offsets4 = @(
    @(54,70)
)
Split-RunsInRange $d $uAbs $offsets4

# ---------------------------------------------------------------------
# 5) ".2 Spil se belezi..." -> "Prelazi se na funkcionalnost cuvanja spila"
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(133)
$p5Start = $p5.Range.Start
$p5End = $p5.Range.End - 1
$target5 = $d.Range($p5Start, $p5End)
$newText5 = ".2 Prelazi se na funkcionalnost čuvanja špila"
$target5.Text = $newText5

$offsets5 = @(
    @(3,3)
)
# split point right after ".2 " (offset 3) - use a 1-length probe just
# after the boundary so Bold toggling affects only the following run
$r5 = $d.Range($p5Start + 3, $p5Start + 4)
$r5.Bold = 1
$r5.Bold = 0

# ---------------------------------------------------------------------
# 6) "Nema." (Preduslovi) -> "Korisnik je ulogovan."
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(144)
$p6Start = $p6.Range.Start
$p6End = $p6.Range.End - 1
$target6 = $d.Range($p6Start, $p6End)
$target6.Text = "Korisnik je ulogovan."

$offsets6 = @(
    @(20,21)
)
Split-RunsInRange $d $p6Start $offsets6

Write-Output "Done"
